$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the Area / Atotal columns (G:H) and the summary columns (J:K)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Column G: per-segment cross-sectional area
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# H2: total area, mirrored into J2/K2 alongside the existing Qtotal
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Match the selection recorded in the saved workbook
$excel.Goto($ws.Range("J2:K2"))
